# Applies the "Save all data to excel" commit to the standards workbook.
# - Row-1 headers on every sheet move from numeric placeholders (1,2[,3,4,5])
#   to real text labels ("Key"/"Value" or "Type"/"Value"/"Unit Rate"/"Amount"/"Notes").
# - A batch of fake/seed data values are refreshed with new random-looking content.
# - "Created at" / "Updated at" rows switch from inline ISO date strings to real
#   Excel date serial numbers formatted as dates.
# - The "Value"-style columns widen slightly to fit the new header text.

$wb = $excel.ActiveWorkbook

$DateFormat = "mm-dd-yy"

function Set-KeyValueHeader($ws) {
    $ws.Range("A1").Value = "Key"
    $ws.Range("B1").Value = "Value"
}

function Set-CreatedUpdated($ws, $createdCell, $createdValue, $updatedCell, $updatedValue) {
    $ws.Range($createdCell).Value = $createdValue
    $ws.Range($createdCell).NumberFormat = $DateFormat
    $ws.Range($updatedCell).Value = $updatedValue
    $ws.Range($updatedCell).NumberFormat = $DateFormat
}

# ---------------------------------------------------------------------------
# Sheet 1: Supplier Data
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Supplier Data")
Set-KeyValueHeader $ws
$ws.Range("B3").Value = "Adam Rodgers"
$ws.Range("B4").Value = "Omnis minus ut iste "
$ws.Range("B5").Value = 2016
$ws.Range("B6").Value = "Laborum Consectetur"
Set-CreatedUpdated $ws "B7" 45690.64534295795 "B8" 45692.01832649687
$ws.Columns.Item(2).ColumnWidth = 27.885714285714286

# ---------------------------------------------------------------------------
# Sheet 2: Fire Alarm Control Panel
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fire Alarm Control Panel")
Set-KeyValueHeader $ws
$ws.Range("B3").Value = "UL"
$ws.Range("B4").Value = 80
$ws.Range("B5").Value = 601
$ws.Range("B6").Value = 101
$ws.Range("B7").Value = 28
$ws.Range("B8").Value = 56
$ws.Range("B9").Value = 100
$ws.Range("B10").Value = 25
$ws.Range("B11").Value = "Esse blanditiis mol"
$ws.Range("B12").Value = 92
$ws.Range("B13").Value = 48
$ws.Range("B14").Value = 90
$ws.Range("B15").Value = 59
$ws.Range("B16").Value = "Voluptatum dolore si"
$ws.Range("B17").Value = "Velit vel dolor non"
$ws.Range("B18").Value = "Aut voluptas nihil e"
$ws.Range("B19").Value = 45
$ws.Range("B20").Value = 64
Set-CreatedUpdated $ws "B21" 45690.64534335464 "B22" 45692.01832656486
$ws.Columns.Item(2).ColumnWidth = 27.885714285714286

# ---------------------------------------------------------------------------
# Sheet 3: Product Data
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Product Data")
Set-KeyValueHeader $ws
$ws.Range("B3").Value = "Unde labore assumend"
$ws.Range("B4").Value = "Mollitia excepturi v"
$ws.Range("B5").Value = "Rerum sit dolor qua"
$ws.Range("B6").Value = 1989
$ws.Range("B7").Value = "Temporibus do obcaec"
$ws.Range("B8").Value = "Iusto quae irure cil"
$ws.Range("B9").Value = "Eum autem aliquam la"
Set-CreatedUpdated $ws "B11" 45690.645347796984 "B12" 45692.0183270467
$ws.Columns.Item(2).ColumnWidth = 27.885714285714286

# ---------------------------------------------------------------------------
# Sheet 4: Graphic Systems
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Graphic Systems")
Set-KeyValueHeader $ws
$ws.Range("B3").Value = "Nisi quis esse vero"
$ws.Range("B4").Value = "Quae dolor proident"
$ws.Range("B5").Value = "Culpa et totam elig"
$ws.Range("B6").Value = "Labore odit enim qui"
$ws.Range("B7").Value = "Excepturi voluptas q"
Set-CreatedUpdated $ws "B9" 45690.64534860172 "B10" 45692.01832712396
$ws.Columns.Item(2).ColumnWidth = 27.885714285714286

# ---------------------------------------------------------------------------
# Sheet 5: Detectors Field Devices
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Detectors Field Devices")
$ws.Range("A1").Value = "Type"
$ws.Range("B1").Value = "Value"
$ws.Range("C1").Value = "Unit Rate"
$ws.Range("D1").Value = "Amount"
$ws.Range("E1").Value = "Notes"

$rows5 = @(
    @{ r = 2;  b = 82; c = 63; d = 5166; e = "Quis ut est minima " },
    @{ r = 3;  b = 39; c = 94; d = 3666; e = "Lorem sunt sit dolor" },
    @{ r = 4;  b = 84; c = 82; d = 6888; e = "Deserunt suscipit es" },
    @{ r = 5;  b = 65; c = 39; d = 2535; e = "Facilis aut dolores " },
    @{ r = 6;  b = 48; c = 98; d = 4704; e = "Veritatis consectetu" },
    @{ r = 7;  b = 81; c = 23; d = 1863; e = "Laboriosam tempore" },
    @{ r = 8;  b = 74; c = 26; d = 1924; e = "Dolor est aut autem" },
    @{ r = 9;  b = 6;  c = 88; d = 528;  e = "Repudiandae nobis cu" },
    @{ r = 10; b = 4;  c = 11; d = 44;   e = "Deleniti sed et sit " },
    @{ r = 11; b = 17; c = 97; d = 1649; e = "Ea enim ipsam fugiat" },
    @{ r = 12; b = 43; c = 34; d = 1462; e = "Amet quis dolores i" },
    @{ r = 13; b = 85; c = 56; d = 4760; e = "Qui nemo delectus e" },
    @{ r = 14; b = 74; c = 87; d = 6438; e = "Et voluptate nemo Na" },
    @{ r = 15; b = 35; c = 23; d = 805;  e = "Sunt quisquam mollit" },
    @{ r = 16; b = 100;c = 89; d = 8900; e = "Occaecat voluptatem" },
    @{ r = 17; b = 58; c = 23; d = 1334; e = "Reprehenderit cupid" },
    @{ r = 18; b = 34; c = 66; d = 2244; e = "Velit sequi odio mo" },
    @{ r = 19; b = 56; c = 73; d = 4088; e = "Temporibus qui repel" },
    @{ r = 20; b = 91; c = 3;  d = 273;  e = "Sequi nulla ad volup" },
    @{ r = 21; b = 67; c = 41; d = 2747; e = "Deserunt et dolorem " }
)
foreach ($row in $rows5) {
    $ws.Range("B" + $row.r).Value = $row.b
    $ws.Range("C" + $row.r).Value = $row.c
    $ws.Range("D" + $row.r).Value = $row.d
    $ws.Range("E" + $row.r).Value = $row.e
}
$ws.Columns.Item(2).ColumnWidth = 8.857142857142858
$ws.Columns.Item(3).ColumnWidth = 13.142857142857142
$ws.Columns.Item(4).ColumnWidth = 9.857142857142858

# ---------------------------------------------------------------------------
# Sheet 6: Manual Pull Station
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Manual Pull Station")
Set-KeyValueHeader $ws
$ws.Range("B3").Value = "Ut aperiam esse dict"
Set-CreatedUpdated $ws "B7" 45690.645346053454 "B8" 45692.01832689102
$ws.Columns.Item(2).ColumnWidth = 27.885714285714286

# ---------------------------------------------------------------------------
# Sheet 7: Door Holders
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Door Holders")
$ws.Range("A1").Value = "Type"
$ws.Range("B1").Value = "Value"
$ws.Range("C1").Value = "Unit Rate"
$ws.Range("D1").Value = "Amount"
$ws.Range("E1").Value = "Notes"

$ws.Range("B2").Value = 80
$ws.Range("C2").Value = 78
$ws.Range("D2").Value = 6240
$ws.Range("E2").Value = "Sint saepe aliqua "

$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 56
$ws.Range("D3").Value = 168
$ws.Range("E3").Value = "Dolor sit non ratio"

$ws.Columns.Item(2).ColumnWidth = 8.857142857142858
$ws.Columns.Item(3).ColumnWidth = 13.142857142857142
$ws.Columns.Item(4).ColumnWidth = 9.857142857142858
$ws.Columns.Item(5).ColumnWidth = 24.142857142857142

# ---------------------------------------------------------------------------
# Sheet 8: Notification Devices
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Notification Devices")
Set-KeyValueHeader $ws
$ws.Range("B3").Value = "Amet illo ea eum au"
$ws.Range("B4").Value = 93
$ws.Range("B6").Value = 5
$ws.Range("B7").Value = 29
$ws.Range("B8").Value = 53
$ws.Range("B9").Value = 65
Set-CreatedUpdated $ws "B11" 45690.645349403334 "B12" 45692.018326822916
$ws.Columns.Item(2).ColumnWidth = 27.885714285714286

Write-Output "edit complete"
